# Swap odds/result data between pairs of rows (keep column A / index untouched)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(13, 14),
    @(20, 21),
    @(22, 23),
    @(26, 27),
    @(28, 29),
    @(30, 31),
    @(51, 52),
    @(54, 55)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value2 = $v2
    $range2.Value2 = $v1
}
